# Rename the embedded logo pictures:
#   - Header: BTec_Logo-Orange   image1.jpg -> image2.jpg
#   - Footer 1: PearsonLogo.png  image2.png -> image1.png
#   - Footer 2: PearsonLogo.png  image2.png -> image1.png
#
# The images live in headers/footers, so we reach them through
# Sections(1).Headers / .Footers rather than Document.InlineShapes
# (which only covers the main story).

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Header: BTec_Logo-Orange -> image2.jpg ---
for ($i = 1; $i -le 3; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
            $shp = $hdr.Range.InlineShapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}

# --- Footers: PearsonLogo.png -> image1.png ---
for ($i = 1; $i -le 3; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
            $shp = $ftr.Range.InlineShapes.Item($j)
            if ($shp.AlternativeText -like "*PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
